# Lesson 1.1 - "added videos to lessons"
#
# A new "Video Time" slide (Title + Content placeholder with a link to a
# YouTube video) is inserted right before the final "Practice" slide, so it
# becomes the new slide 8 and the old "Practice" slide is pushed down to
# become slide 9.

$p = $ppt.ActivePresentation

# The existing "Practice" slide is slide 8 (last slide). Insert the new
# slide in front of it using the same "Title and Content" layout (layout
# index 2, same as the Practice slide uses).
$s = $p.Slides.Add(8, 2)

# Title placeholder.
$title = $s.Shapes.Placeholders.Item(1)
$title.TextFrame.TextRange.Text = "Video Time"

# Content placeholder: a hyperlinked line with the video URL followed by a
# blank line.
$body = $s.Shapes.Placeholders.Item(2)
$url = "https://www.youtube.com/watch?v=gpH8T2CRlLI"

$bodyRange = $body.TextFrame.TextRange
$bodyRange.Text = $url
$bodyRange.ActionSettings.Item(1).Action = 7
$bodyRange.ActionSettings.Item(1).Hyperlink.Address = $url

# Append a second, empty paragraph after the hyperlinked line.
$current = $body.TextFrame.TextRange.Text
$body.TextFrame.TextRange.Text = $current + "`rX"
$secondPara = $body.TextFrame.TextRange.Paragraphs(2, 1)
$secondPara.Text = ""
